$wb = $excel.ActiveWorkbook

# Work on the "Nk" worksheet (sheet2.xml): add the new rows of data
$ws = $wb.Worksheets.Item("Nk")

$ws.Range("A10").Value = " antilope"
$ws.Range("A11").Value = "Antilope"
$ws.Range("A12").Value = "Barbagianni"
$ws.Range("A13").Value = "Pippo"
$ws.Range("A14").Value = "Fagiolo"
$ws.Range("A15").Value = "Aaaaaaaaaaa"
$ws.Range("A16").Value = "Aaaaaaaaaaaaaaaaaa"

# Page setup for the "Nk" sheet
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# Make "Nk" the active sheet / tab, and set its selection
$ws.Activate() | Out-Null
$ws.Range("F10").Select() | Out-Null
